$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) "总计" (summary) sheet: insert a new first data row for "2022-Q3" and
#    shift the existing quarters down by one row, adding "2020-Q4" at the end.
# ---------------------------------------------------------------------------
$summaryWs = $wb.Worksheets.Item("总计")

# Insert a blank row right under the header, pushing all quarter rows down one.
$summaryWs.Rows.Item(2).Insert()
$summaryWs.Range("B2:D2").ClearFormats()

# The row-index column (A) is a simple running counter (0,1,2,...) tied to the
# row position, not to the quarter - after the insert it needs to be bumped
# by one for every pre-existing data row (old row2 -> row3 keeps counter 0,
# but the target wants 1, etc).
for ($r = 9; $r -ge 3; $r--) {
    $cur = $summaryWs.Cells.Item($r, 1).Value2
    $summaryWs.Cells.Item($r, 1).Value = $cur + 1
}

# New first data row: 2022-Q3
$summaryWs.Range("A2").Value = 0
$summaryWs.Range("A3").Copy()
$summaryWs.Range("A2").PasteSpecial(-4122)
$summaryWs.Range("B2").Value = "2022-Q3"
$summaryWs.Range("C2").Value = 2
$summaryWs.Range("D2").Value = 0.09

# ---------------------------------------------------------------------------
# 2) Add the new "2022-Q3" worksheet (fund holder detail), right after "总计".
# ---------------------------------------------------------------------------
$refWs = $wb.Worksheets.Item("2022-Q2")
$newWs = $wb.Worksheets.Add($null, $summaryWs)
$newWs.Name = "2022-Q3"
$newWs.Outline.SummaryRow = 1
$newWs.Outline.SummaryColumn = 1

# Match the header / index-column formatting used by the other quarter sheets.
$refWs.Range("B1:H1").Copy()
$newWs.Range("B1:H1").PasteSpecial(-4122)
$refWs.Range("A2:A3").Copy()
$newWs.Range("A2:A3").PasteSpecial(-4122)

$newWs.Range("B1").Value = "基金代码"
$newWs.Range("C1").Value = "基金名称"
$newWs.Range("D1").Value = "基金规模"
$newWs.Range("E1").Value = "股票总仓位"
$newWs.Range("F1").Value = "仓位占比"
$newWs.Range("G1").Value = "持有市值(亿元)"
$newWs.Range("H1").Value = "仓位排名"

# Keep the numeric-looking strings (fund code / scale / position %) as text,
# matching the source data's inlineStr cell type.
$newWs.Range("B2:G3").NumberFormat = "@"

$newWs.Range("A2").Value = 0
$newWs.Range("B2").Value = "002707"
$newWs.Range("C2").Value = "摩根士丹利华鑫科技领先灵活配置混合A"
$newWs.Range("D2").Value = "1.75"
$newWs.Range("E2").Value = "94.13"
$newWs.Range("F2").Value = "5.07"
$newWs.Range("G2").Value = "0.0887"
$newWs.Range("H2").Value = 2

$newWs.Range("A3").Value = 1
$newWs.Range("B3").Value = "014871"
$newWs.Range("C3").Value = "摩根士丹利华鑫科技领先灵活配置混合C"
$newWs.Range("D3").Value = "0.08"
$newWs.Range("E3").Value = "94.13"
$newWs.Range("F3").Value = "5.07"
$newWs.Range("G3").Value = "0.0041"
$newWs.Range("H3").Value = 2

# Drop the temporary text-number-format so the data cells carry no explicit
# style, matching the other quarter sheets.
$newWs.Range("B2:G3").ClearFormats()

# Restore the originally active sheet/selection.
$summaryWs.Activate() | Out-Null
$summaryWs.Range("A1").Select() | Out-Null

